$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.425.16'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +7.18%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.72'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +6.47%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.43%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.56%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3834'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.51%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.37'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +5.59%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3520'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.233'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.89%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07792'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.40%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.54%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.49'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +12.95%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.636'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.236'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.61%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.813.29'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.65%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001127'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06768'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.86%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.30'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +8.16%  '

# Row 20
$ws.Range("E20").Value = '  +0.47%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.90'
$ws.Range("D21").ClearFormats()

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.559'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +8.67%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.13'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.46%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.410.99'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +7.17%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.474'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.94%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.684'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.59%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.07'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +15.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.514'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +18.69%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.48'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.51%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.018.52'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.81%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '137.21'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +7.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.460'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +8.52%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.167'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.78%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '14.09'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +11.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08793'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.99%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.729'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.50%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.678'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +7.52%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6948'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +14.54%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06554'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.85%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2274'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +7.52%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02421'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +7.07%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.028'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.77%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.254'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.90%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.91'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.57%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6539'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +12.24%  '

# Row 46
$ws.Range("E46").Value = '  +0.48%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.036'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +5.47%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.183'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +9.50%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.24'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.00%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07352'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.27%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.80'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.53%  '
